$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row 173 that is an exact copy (values + formatting) of the
# current row 172. This becomes the final row of the block once the other
# rows are shifted below.
$ws.Rows.Item(172).Copy()
$ws.Rows.Item(173).Insert()

# Shift the "variable" columns (date, volume, prices, unit, region,
# avg-price, kg) of rows 172..167 down by one, taking the values that used
# to sit in the row immediately above each of them. Columns A,B,C,E,F,G,H,I,R
# are identical across this block of rows, so they do not need to change.
$ws.Range("D172").Value  = $ws.Range("D171").Value()
$ws.Range("J172").Value  = $ws.Range("J171").Value()
$ws.Range("K172").Value  = $ws.Range("K171").Value()
$ws.Range("L172").Value  = $ws.Range("L171").Value()
$ws.Range("M172").Value  = $ws.Range("M171").Value()
$ws.Range("N172").Value  = $ws.Range("N171").Value()
$ws.Range("O172").Value  = $ws.Range("O171").Value()
$ws.Range("P172").Value  = $ws.Range("P171").Value()
$ws.Range("Q172").Value  = $ws.Range("Q171").Value()

$ws.Range("D171").Value  = $ws.Range("D170").Value()
$ws.Range("J171").Value  = $ws.Range("J170").Value()
$ws.Range("K171").Value  = $ws.Range("K170").Value()
$ws.Range("L171").Value  = $ws.Range("L170").Value()
$ws.Range("M171").Value  = $ws.Range("M170").Value()
$ws.Range("N171").Value  = $ws.Range("N170").Value()
$ws.Range("O171").Value  = $ws.Range("O170").Value()
$ws.Range("P171").Value  = $ws.Range("P170").Value()
$ws.Range("Q171").Value  = $ws.Range("Q170").Value()

$ws.Range("D170").Value  = $ws.Range("D169").Value()
$ws.Range("J170").Value  = $ws.Range("J169").Value()
$ws.Range("K170").Value  = $ws.Range("K169").Value()
$ws.Range("L170").Value  = $ws.Range("L169").Value()
$ws.Range("M170").Value  = $ws.Range("M169").Value()
$ws.Range("N170").Value  = $ws.Range("N169").Value()
$ws.Range("O170").Value  = $ws.Range("O169").Value()
$ws.Range("P170").Value  = $ws.Range("P169").Value()
$ws.Range("Q170").Value  = $ws.Range("Q169").Value()

$ws.Range("D169").Value  = $ws.Range("D168").Value()
$ws.Range("J169").Value  = $ws.Range("J168").Value()
$ws.Range("K169").Value  = $ws.Range("K168").Value()
$ws.Range("L169").Value  = $ws.Range("L168").Value()
$ws.Range("M169").Value  = $ws.Range("M168").Value()
$ws.Range("N169").Value  = $ws.Range("N168").Value()
$ws.Range("O169").Value  = $ws.Range("O168").Value()
$ws.Range("P169").Value  = $ws.Range("P168").Value()
$ws.Range("Q169").Value  = $ws.Range("Q168").Value()

$ws.Range("D168").Value  = $ws.Range("D167").Value()
$ws.Range("J168").Value  = $ws.Range("J167").Value()
$ws.Range("K168").Value  = $ws.Range("K167").Value()
$ws.Range("L168").Value  = $ws.Range("L167").Value()
$ws.Range("M168").Value  = $ws.Range("M167").Value()
$ws.Range("N168").Value  = $ws.Range("N167").Value()
$ws.Range("O168").Value  = $ws.Range("O167").Value()
$ws.Range("P168").Value  = $ws.Range("P167").Value()
$ws.Range("Q168").Value  = $ws.Range("Q167").Value()

$ws.Range("D167").Value  = $ws.Range("D166").Value()
$ws.Range("J167").Value  = $ws.Range("J166").Value()
$ws.Range("K167").Value  = $ws.Range("K166").Value()
$ws.Range("L167").Value  = $ws.Range("L166").Value()
$ws.Range("M167").Value  = $ws.Range("M166").Value()
$ws.Range("N167").Value  = $ws.Range("N166").Value()
$ws.Range("O167").Value  = $ws.Range("O166").Value()
$ws.Range("P167").Value  = $ws.Range("P166").Value()
$ws.Range("Q167").Value  = $ws.Range("Q166").Value()

# Row 166 gets a brand-new record: new report date and new volume, the
# rest of its fields (price, unit, region, avg price, kg) are unchanged.
$ws.Range("D166").Value = 44509
$ws.Range("J166").Value = 160
